# increase aluop signal to 3-bits
# - widen the ALUOp (column J) control codes from 2-bit to 3-bit encodings
# - add a new "sign" control column (N) for every instruction row

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New column header
$ws.Range("N1").Value = 'sign'

# add
$ws.Range("J2").Value = "2'b010"
$ws.Range("N2").Value = 'x'

# sub
$ws.Range("J3").Value = "2'b010"
$ws.Range("N3").Value = 'x'

# mul
$ws.Range("J4").Value = "2'b010"
$ws.Range("N4").Value = 'x'

# and
$ws.Range("J5").Value = "2'b010"
$ws.Range("N5").Value = 'x'

# or
$ws.Range("J6").Value = "2'b010"
$ws.Range("N6").Value = 'x'

# nor
$ws.Range("J7").Value = "2'b010"
$ws.Range("N7").Value = 'x'

# slt
$ws.Range("J8").Value = "2'b010"
$ws.Range("N8").Value = 'x'

# sltu
$ws.Range("N9").Value = 'x'

# sll
$ws.Range("J10").Value = "2'b010"
$ws.Range("N10").Value = 'x'

# srl
$ws.Range("N11").Value = 'x'

# jr
$ws.Range("J12").Value = "2'bxxx"
$ws.Range("N12").Value = 'x'

# slti
$ws.Range("N13").Value = 1

# sltui
$ws.Range("N14").Value = 0

# addi
$ws.Range("J15").Value = "2'b000"
$ws.Range("N15").Value = 1

# andi
$ws.Range("J16").Value = "2'b011"
$ws.Range("N16").Value = 0

# ori (previously only had its instruction label; now fully populated)
$ws.Range("B17").Value = "b'001101"
$ws.Range("C17").Value = "2'b00"
$ws.Range("D17").Value = 1
$ws.Range("E17").Value = "2'b00"
$ws.Range("F17").Value = 1
$ws.Range("G17").Value = 0
$ws.Range("H17").Value = 0
$ws.Range("I17").Value = 0
$ws.Range("J17").Value = "2'b011"
$ws.Range("K17").Value = 0
$ws.Range("L17").Value = "2'bxx"
$ws.Range("M17").Value = 'x'
$ws.Range("N17").Value = 0

# lui
$ws.Range("N18").Value = 0

# lw
$ws.Range("J19").Value = "2'b000"
$ws.Range("N19").Value = 1

# lh
$ws.Range("J20").Value = "2'b000"
$ws.Range("N20").Value = 1

# lb
$ws.Range("J21").Value = "2'b000"
$ws.Range("N21").Value = 1

# lhu
$ws.Range("J22").Value = "2'b000"
$ws.Range("N22").Value = 1

# lbu
$ws.Range("J23").Value = "2'b000"
$ws.Range("N23").Value = 1

# sw
$ws.Range("J24").Value = "2'b000"
$ws.Range("N24").Value = 1

# sh
$ws.Range("J25").Value = "2'b000"
$ws.Range("N25").Value = 1

# sb
$ws.Range("J26").Value = "2'b000"
$ws.Range("N26").Value = 1

# beq
$ws.Range("J27").Value = "2'b001"
$ws.Range("N27").Value = 1

# bne
$ws.Range("N28").Value = 1

# j
$ws.Range("J29").Value = "2'bxxx"
$ws.Range("N29").Value = 'x'

# jal
$ws.Range("J30").Value = "2'bxxx"
$ws.Range("N30").Value = 'x'

# leave the active selection where the author ended up
$ws.Range("J31").Select()
